$wb = $excel.ActiveWorkbook

# zh-cn sheet: row for 484b79dc... (row 3) holds the Handoff/Handback datetimes
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-24 02:24:15"
$wsZhCn.Range("H3").Value = "2016-03-24 02:25:04"

# de-de sheet: row for 484b79dc... (row 3) holds the Handoff/Handback datetimes
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-24 02:24:24"
$wsDeDe.Range("H3").Value = "2016-03-24 02:25:18"
